$wb = $excel.ActiveWorkbook

# --- Sheet "展览" (Exhibition) ---
$ws = $wb.Worksheets.Item("展览")
$ws.Range("F4").Value = 1186
$ws.Range("F6").Value = 2797
$ws.Range("I6").Value = "//i2.hdslb.com/bfs/openplatform/202410/hHMKXwCN1729135895623.jpeg"
$ws.Range("F9").Value = 135
$ws.Range("F10").Value = 303
$ws.Range("F12").Value = 717
$ws.Range("F13").Value = 122
$ws.Range("F15").Value = 1888
$ws.Range("F17").Value = 47
$ws.Range("F18").Value = 210

# --- Sheet "演出" (Show) ---
$ws = $wb.Worksheets.Item("演出")
$ws.Range("F7").Value = 29
$ws.Range("F10").Value = 70
$ws.Range("F13").Value = 217

# --- Sheet "本地生活" (Local Life) ---
$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F2").Value = 6375
$ws.Range("F3").Value = 802
$ws.Range("F4").Value = 2045
$ws.Range("F5").Value = 279

# --- Sheet "全部类型" (All Types) ---
$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F2").Value = 6375
$ws.Range("F3").Value = 802
$ws.Range("F4").Value = 2045
$ws.Range("F5").Value = 279
$ws.Range("F12").Value = 1186
$ws.Range("F15").Value = 29
$ws.Range("F17").Value = 2797
$ws.Range("I17").Value = "//i2.hdslb.com/bfs/openplatform/202410/hHMKXwCN1729135895623.jpeg"
$ws.Range("F20").Value = 70
$ws.Range("F24").Value = 135
$ws.Range("F25").Value = 303
$ws.Range("F26").Value = 217
$ws.Range("F28").Value = 717
$ws.Range("F29").Value = 122
$ws.Range("F32").Value = 1888
$ws.Range("F36").Value = 47
$ws.Range("F37").Value = 210
